$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$h = $sec.Headers.Item(1)
$rng = $h.Range
$found = $rng.Find.Execute("image2.jpg", $true, $false, $false, $false, $false, $true, 1, $false, "image1.jpg", 2)
Write-Host "Found: $found"
